$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.947.75'
$ws.Range('E2').Value = '  -3.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.912.37'
$ws.Range('E3').Value = '  -4.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.80'
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.48'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -3.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.908.29'
$ws.Range('E9').Value = '  -4.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.73'
$ws.Range('E10').Value = '  +5.74%  '
$ws.Range('E11').Value = '  -5.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -2.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000225'
$ws.Range('E13').Value = '  -4.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.52'
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.396.07'
$ws.Range('E16').Value = '  -4.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.862.15'
$ws.Range('E17').Value = '  -3.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.77'
$ws.Range('E18').Value = '  -3.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.912.36'
$ws.Range('E19').Value = '  -4.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '427.25'
$ws.Range('E20').Value = '  -6.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.58'
$ws.Range('E21').Value = '  -4.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.672'
$ws.Range('E22').Value = '  -3.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.12'
$ws.Range('E23').Value = '  -5.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.10'
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.96'
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('E26').Value = '  -2.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.84'
$ws.Range('E27').Value = '  -3.05%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.26'
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.62'
$ws.Range('E31').Value = '  -3.55%  '
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.44'
$ws.Range('E33').Value = '  -4.48%  '
$ws.Range('E34').Value = '  -3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0847'
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('E36').Value = '  -2.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.63'
$ws.Range('E37').Value = '  -4.69%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.98'
$ws.Range('E38').Value = '  -4.76%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.41'
$ws.Range('E39').Value = '  -2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').Value = '  -4.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('E42').Value = '  -4.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.292'
$ws.Range('E43').Value = '  -0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.66'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '376.87'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0347'
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.673.35'
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.81'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.83'
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('E51').Value = '  -1.79%  '
